$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new data row (A8:F8) matching the existing table's schema:
# unitTypeId, projectId, name, sellingPrice, available, total
$ws.Range("A8").Value = 4
$ws.Range("B8").Value = 3
$ws.Range("C8").Value = "2-Room"
$ws.Range("D8").Value = 5000
$ws.Range("E8").Value = 50
$ws.Range("F8").Value = 50
